$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 9906
$ws.Range("J7").Value = 9906
$ws.Range("L7").Value = 9906
$ws.Range("N7").Value = -10130
$ws.Range("H14").Value = 9906
$ws.Range("J14").Value = 9906
$ws.Range("L14").Value = 9906
$ws.Range("N14").Value = -10288
$ws.Range("H17").Value = 1616.1333
$ws.Range("J17").Value = 1616.1333
$ws.Range("L17").Value = 4848.3999
$ws.Range("N17").Value = -5184.3999
$ws.Range("H53").Value = 472.42856
$ws.Range("I53").Value = 304
$ws.Range("K53").Value = 304
$ws.Range("M53").Value = 333
$ws.Range("H64").Value = 4822
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 4822
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H98").Value = 963.625
$ws.Range("I98").Value = 963.625
$ws.Range("K98").Value = 963.625
$ws.Range("M98").Value = 534.375
$ws.Range("H107").Value = 582.1667
$ws.Range("I107").Value = 498.5
$ws.Range("J107").Value = 749.5
$ws.Range("K107").Value = 498.5
$ws.Range("L107").Value = 749.5
$ws.Range("M107").Value = 1421.5
$ws.Range("N107").Value = -4589.5
$ws.Range("H112").Value = 1878.5862
$ws.Range("J112").Value = 1999.5834
$ws.Range("L112").Value = 5998.7502
$ws.Range("N112").Value = -8214.7502
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H120").Value = 97000
$ws.Range("J120").Value = 97000
$ws.Range("L120").Value = 97000
$ws.Range("N120").Value = -106676
$ws.Range("H122").Value = 963.625
$ws.Range("I122").Value = 963.625
$ws.Range("K122").Value = 2890.875
$ws.Range("M122").Value = -440.875
$ws.Range("H123").Value = 66333.336
$ws.Range("J123").Value = 66333.336
$ws.Range("L123").Value = 66333.336
$ws.Range("N123").Value = -76133.336
$ws.Range("H124").Value = 67250
$ws.Range("J124").Value = 67250
$ws.Range("L124").Value = 67250
$ws.Range("N124").Value = -77070
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 5815.9473
$ws.Range("I137").Value = 3468.875
$ws.Range("K137").Value = 10406.625
$ws.Range("M137").Value = -7856.625
$ws.Range("H138").Value = 5477.148
$ws.Range("J138").Value = 5570.25
$ws.Range("L138").Value = 16710.75
$ws.Range("N138").Value = -26990.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3585
$ws.Range("I61").Value = 2879.2
$ws.Range("J61").Value = 4026.125
$ws.Range("K61").Value = 2879.2
$ws.Range("L61").Value = 4026.125
$ws.Range("M61").Value = -2667.2
$ws.Range("N61").Value = -4450.125
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H132").Value = 3780.4546
$ws.Range("I132").Value = 2211.5715
$ws.Range("K132").Value = 6634.7145
$ws.Range("M132").Value = -4104.7145
$ws.Range("H136").Value = 3585
$ws.Range("I136").Value = 2879.2
$ws.Range("J136").Value = 4026.125
$ws.Range("K136").Value = 8637.599999999999
$ws.Range("L136").Value = 12078.375
$ws.Range("M136").Value = -6087.599999999999
$ws.Range("N136").Value = -17178.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1581
$ws.Range("I5").Value = 1581
$ws.Range("K5").Value = 1581
$ws.Range("M5").Value = -1468
$ws.Range("H134").Value = 2817.0667
$ws.Range("I134").Value = 1778.7778
$ws.Range("K134").Value = 5336.3334
$ws.Range("M134").Value = -2801.3334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6690.6772
$ws.Range("I31").Value = 3731.682
$ws.Range("K31").Value = 3731.682
$ws.Range("M31").Value = -3436.682
$ws.Range("H34").Value = 6690.6772
$ws.Range("I34").Value = 3731.682
$ws.Range("K34").Value = 3731.682
$ws.Range("M34").Value = -3529.682
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H122").Value = 4882.4
$ws.Range("I122").Value = 4882.4
$ws.Range("K122").Value = 14647.2
$ws.Range("M122").Value = -12197.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1245.9
$ws.Range("J122").Value = 1188.1428
$ws.Range("L122").Value = 10693.2852
$ws.Range("N122").Value = -15593.2852
$ws.Range("H131").Value = 46144.24
$ws.Range("J131").Value = 1943.7778
$ws.Range("L131").Value = 5831.3334
$ws.Range("N131").Value = -15911.3334

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3799
$ws.Range("J122").Value = 2999.5
$ws.Range("L122").Value = 8998.5
$ws.Range("N122").Value = -13898.5
$ws.Range("H132").Value = 1669.875
$ws.Range("I132").Value = 1572
$ws.Range("K132").Value = 4716
$ws.Range("M132").Value = -2186

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1640
$ws.Range("I22").Value = 1600
$ws.Range("J22").Value = 1650
$ws.Range("K22").Value = 1600
$ws.Range("L22").Value = 1650
$ws.Range("M22").Value = -1305
$ws.Range("N22").Value = -2240
$ws.Range("H27").Value = 1640
$ws.Range("I27").Value = 1600
$ws.Range("J27").Value = 1650
$ws.Range("K27").Value = 1600
$ws.Range("L27").Value = 1650
$ws.Range("M27").Value = -1493
$ws.Range("N27").Value = -1864
$ws.Range("H40").Value = 4303.6665
$ws.Range("I40").Value = 3999
$ws.Range("J40").Value = 4390.7144
$ws.Range("K40").Value = 3999
$ws.Range("L40").Value = 4390.7144
$ws.Range("M40").Value = -3863
$ws.Range("N40").Value = -4662.7144
$ws.Range("H55").Value = 636.8
$ws.Range("J55").Value = 427.66666
$ws.Range("L55").Value = 427.66666
$ws.Range("N55").Value = -773.66666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6489.0713
$ws.Range("I122").Value = 5723.5454
$ws.Range("K122").Value = 17170.6362
$ws.Range("M122").Value = -14720.6362
